$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1) to reflect new column labels.
# Unchanged cells (A1, B1, C1, D1, H1, J1) are left as-is.
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

# Update the active selection to the header row range instead of a single cell.
$ws.Range("A1:K1").Select()
